# Generate Report for Handback
#
# Previously-outstanding files (8e14a4a7-...md and d01e54dd-...md) have now
# been handed back (status flips from "Ready for handoff" to "Handed back:
# in sync with en-US" everywhere), rows are re-sorted alphabetically by file
# name (237dd7c3 < 8e14a4a7 < aabe33c5 < d01e54dd), and the "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns are
# now populated for every row on the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$statusHanded = "Handed back: in sync with en-US"
$extMd = ".md"
$include = "Include"

# ======================================================================
# Overview sheet
# ======================================================================
$ov = $wb.Worksheets.Item("Overview")

# Row 2 (237dd7c3) is already correct/unchanged.

# Row 3 now holds 8e14a4a7 (previously row 4) -- now handed back.
$ov.Range("A3").Value = "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md"
$ov.Range("B3").Value = $statusHanded
$ov.Range("C3").Value = $statusHanded
$ov.Range("D3").Value = "2016-03-22 06:22:45"

# Row 4 now holds aabe33c5 (previously row 3).
$ov.Range("A4").Value = "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md"
$ov.Range("B4").Value = $statusHanded
$ov.Range("C4").Value = $statusHanded
$ov.Range("D4").Value = "2016-03-22 06:21:43"

# Row 5 (d01e54dd) stays last -- now handed back.
$ov.Range("A5").Value = "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md"
$ov.Range("B5").Value = $statusHanded
$ov.Range("C5").Value = $statusHanded
$ov.Range("D5").Value = "2016-03-22 06:22:45"

# Rebuild the hyperlinks on column A in the new row order.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4bf54648002af0c28612189322737d3e1f0f8d50/e2e/237dd7c3-c39c-4765-b965-031b913e1a2e.md", "", "", "237dd7c3-c39c-4765-b965-031b913e1a2e.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/054150cb1a65e91fcf4cef33d21f7c2ffd5a32d1/e2e/8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md", "", "", "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4bf54648002af0c28612189322737d3e1f0f8d50/e2e/aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md", "", "", "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md")
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/054150cb1a65e91fcf4cef33d21f7c2ffd5a32d1/e2e/d01e54dd-f5f1-4163-97f9-d5a2625eda64.md", "", "", "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md")

# ======================================================================
# zh-cn sheet
# ======================================================================
$zh = $wb.Worksheets.Item("zh-cn")

# --- Row 2 : 237dd7c3 (unchanged values, kept explicit for clarity) ---
$zh.Range("A2").Value = "237dd7c3-c39c-4765-b965-031b913e1a2e.md"
$zh.Range("B2").Value = $extMd
$zh.Range("C2").Value = $statusHanded
$zh.Range("D2").Value = "237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-22 06:21:38"
$zh.Range("F2").Value = "237dd7c3-c39c-4765-b965-031b913e1a2e.md"
$zh.Range("G2").Value = "237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-22 06:22:04"
$zh.Range("J2").Value = $include

# --- Row 3 : 8e14a4a7 (now handed back, moved up from old row 4) ---
$zh.Range("A3").Value = "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md"
$zh.Range("B3").Value = $extMd
$zh.Range("C3").Value = $statusHanded
$zh.Range("D3").Value = "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-22 06:22:41"
$zh.Range("F3").Value = "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md"
$zh.Range("G3").Value = "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-22 06:23:06"
$zh.Range("J3").Value = $include

# --- Row 4 : aabe33c5 (moved down from old row 3) ---
$zh.Range("A4").Value = "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md"
$zh.Range("B4").Value = $extMd
$zh.Range("C4").Value = $statusHanded
$zh.Range("D4").Value = "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.zh-cn.xlf"
$zh.Range("E4").Value = "2016-03-22 06:21:38"
$zh.Range("F4").Value = "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md"
$zh.Range("G4").Value = "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.zh-cn.xlf"
$zh.Range("H4").Value = "2016-03-22 06:22:04"
$zh.Range("J4").Value = $include

# --- Row 5 : d01e54dd (now handed back, stays last) ---
$zh.Range("A5").Value = "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md"
$zh.Range("B5").Value = $extMd
$zh.Range("C5").Value = $statusHanded
$zh.Range("D5").Value = "d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.zh-cn.xlf"
$zh.Range("E5").Value = "2016-03-22 06:22:41"
$zh.Range("F5").Value = "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md"
$zh.Range("G5").Value = "d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.zh-cn.xlf"
$zh.Range("H5").Value = "2016-03-22 06:23:06"
$zh.Range("J5").Value = $include

# Ensure the newly-populated F/G columns render with the same (hyperlink)
# style used by the rest of column A/D/F/G.
$zh.Range("F3:G5").Style = $zh.Range("F2").Style

# Rebuild hyperlinks on A / D / F / G for all four data rows, in row order,
# so relationship ids come out sequential.
$zh.Hyperlinks.Delete()

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4bf54648002af0c28612189322737d3e1f0f8d50/e2e/237dd7c3-c39c-4765-b965-031b913e1a2e.md", "", "", "237dd7c3-c39c-4765-b965-031b913e1a2e.md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ca9fe4f069ae3b0d90a3c8c6e7f13bfb7eef0d3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.zh-cn.xlf", "", "", "237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c1fde1770ae00a8122197e45ea76b4dd6b97f4e6/e2e/237dd7c3-c39c-4765-b965-031b913e1a2e.md", "", "", "237dd7c3-c39c-4765-b965-031b913e1a2e.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b6c4727899f4116b7946ac8887e355c59d8aa763/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.zh-cn.xlf", "", "", "237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/054150cb1a65e91fcf4cef33d21f7c2ffd5a32d1/e2e/8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md", "", "", "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c1950be0630afd63f695f1ca492765aa7ed1bafa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.zh-cn.xlf", "", "", "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c1fde1770ae00a8122197e45ea76b4dd6b97f4e6/e2e/8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md", "", "", "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md")
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b6c4727899f4116b7946ac8887e355c59d8aa763/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.zh-cn.xlf", "", "", "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4bf54648002af0c28612189322737d3e1f0f8d50/e2e/aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md", "", "", "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md")
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ca9fe4f069ae3b0d90a3c8c6e7f13bfb7eef0d3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.zh-cn.xlf", "", "", "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c1fde1770ae00a8122197e45ea76b4dd6b97f4e6/e2e/aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md", "", "", "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md")
$zh.Hyperlinks.Add($zh.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b6c4727899f4116b7946ac8887e355c59d8aa763/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.zh-cn.xlf", "", "", "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/054150cb1a65e91fcf4cef33d21f7c2ffd5a32d1/e2e/d01e54dd-f5f1-4163-97f9-d5a2625eda64.md", "", "", "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md")
$zh.Hyperlinks.Add($zh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c1950be0630afd63f695f1ca492765aa7ed1bafa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.zh-cn.xlf", "", "", "d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c1fde1770ae00a8122197e45ea76b4dd6b97f4e6/e2e/d01e54dd-f5f1-4163-97f9-d5a2625eda64.md", "", "", "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md")
$zh.Hyperlinks.Add($zh.Range("G5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b6c4727899f4116b7946ac8887e355c59d8aa763/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.zh-cn.xlf", "", "", "d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.zh-cn.xlf")

# ======================================================================
# de-de sheet
# ======================================================================
$de = $wb.Worksheets.Item("de-de")

# --- Row 2 : 237dd7c3 (unchanged values, kept explicit for clarity) ---
$de.Range("A2").Value = "237dd7c3-c39c-4765-b965-031b913e1a2e.md"
$de.Range("B2").Value = $extMd
$de.Range("C2").Value = $statusHanded
$de.Range("D2").Value = "237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.de-de.xlf"
$de.Range("E2").Value = "2016-03-22 06:21:43"
$de.Range("F2").Value = "237dd7c3-c39c-4765-b965-031b913e1a2e.md"
$de.Range("G2").Value = "237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.de-de.xlf"
$de.Range("H2").Value = "2016-03-22 06:22:10"
$de.Range("J2").Value = $include

# --- Row 3 : 8e14a4a7 (now handed back, moved up from old row 4) ---
$de.Range("A3").Value = "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md"
$de.Range("B3").Value = $extMd
$de.Range("C3").Value = $statusHanded
$de.Range("D3").Value = "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.de-de.xlf"
$de.Range("E3").Value = "2016-03-22 06:22:45"
$de.Range("F3").Value = "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md"
$de.Range("G3").Value = "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.de-de.xlf"
$de.Range("H3").Value = "2016-03-22 06:23:12"
$de.Range("J3").Value = $include

# --- Row 4 : aabe33c5 (moved down from old row 3) ---
$de.Range("A4").Value = "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md"
$de.Range("B4").Value = $extMd
$de.Range("C4").Value = $statusHanded
$de.Range("D4").Value = "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.de-de.xlf"
$de.Range("E4").Value = "2016-03-22 06:21:43"
$de.Range("F4").Value = "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md"
$de.Range("G4").Value = "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.de-de.xlf"
$de.Range("H4").Value = "2016-03-22 06:22:10"
$de.Range("J4").Value = $include

# --- Row 5 : d01e54dd (now handed back, stays last) ---
$de.Range("A5").Value = "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md"
$de.Range("B5").Value = $extMd
$de.Range("C5").Value = $statusHanded
$de.Range("D5").Value = "d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.de-de.xlf"
$de.Range("E5").Value = "2016-03-22 06:22:45"
$de.Range("F5").Value = "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md"
$de.Range("G5").Value = "d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.de-de.xlf"
$de.Range("H5").Value = "2016-03-22 06:23:12"
$de.Range("J5").Value = $include

# Ensure the newly-populated F/G columns render with the same (hyperlink)
# style used by the rest of column A/D/F/G.
$de.Range("F3:G5").Style = $de.Range("F2").Style

# Rebuild hyperlinks on A / D / F / G for all four data rows, in row order,
# so relationship ids come out sequential.
$de.Hyperlinks.Delete()

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4bf54648002af0c28612189322737d3e1f0f8d50/e2e/237dd7c3-c39c-4765-b965-031b913e1a2e.md", "", "", "237dd7c3-c39c-4765-b965-031b913e1a2e.md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/686e99a75da1b6dc5b893020a0bfb30f1b6c1f60/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.de-de.xlf", "", "", "237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1f25d0bb1fa5aca76e47d2e4a727b364fa77cdb0/e2e/237dd7c3-c39c-4765-b965-031b913e1a2e.md", "", "", "237dd7c3-c39c-4765-b965-031b913e1a2e.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c6931919e8bacbd713899fa79f7ba0d9aa26f5e1/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.de-de.xlf", "", "", "237dd7c3-c39c-4765-b965-031b913e1a2e.ffde92be40841bafa2c277373ea0e7e4497178ff.de-de.xlf")

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/054150cb1a65e91fcf4cef33d21f7c2ffd5a32d1/e2e/8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md", "", "", "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f0906536d8fa515ff3deee7e59cd578785908c8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.de-de.xlf", "", "", "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1f25d0bb1fa5aca76e47d2e4a727b364fa77cdb0/e2e/8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md", "", "", "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.md")
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c6931919e8bacbd713899fa79f7ba0d9aa26f5e1/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.de-de.xlf", "", "", "8e14a4a7-4766-4ae2-a79c-de85c5dc3078.602c85742be0d0336097d89580e5c901993db122.de-de.xlf")

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4bf54648002af0c28612189322737d3e1f0f8d50/e2e/aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md", "", "", "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md")
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/686e99a75da1b6dc5b893020a0bfb30f1b6c1f60/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.de-de.xlf", "", "", "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1f25d0bb1fa5aca76e47d2e4a727b364fa77cdb0/e2e/aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md", "", "", "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.md")
$de.Hyperlinks.Add($de.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c6931919e8bacbd713899fa79f7ba0d9aa26f5e1/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.de-de.xlf", "", "", "aabe33c5-c4ef-4953-a8ac-3f4185db69f8.572d081cfe1a0846fb700f4ad47b8e6d82c0f08b.de-de.xlf")

$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/054150cb1a65e91fcf4cef33d21f7c2ffd5a32d1/e2e/d01e54dd-f5f1-4163-97f9-d5a2625eda64.md", "", "", "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md")
$de.Hyperlinks.Add($de.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f0906536d8fa515ff3deee7e59cd578785908c8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.de-de.xlf", "", "", "d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1f25d0bb1fa5aca76e47d2e4a727b364fa77cdb0/e2e/d01e54dd-f5f1-4163-97f9-d5a2625eda64.md", "", "", "d01e54dd-f5f1-4163-97f9-d5a2625eda64.md")
$de.Hyperlinks.Add($de.Range("G5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c6931919e8bacbd713899fa79f7ba0d9aa26f5e1/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.de-de.xlf", "", "", "d01e54dd-f5f1-4163-97f9-d5a2625eda64.357abd4dcd8e664313864a55bdf2f32ce86fb947.de-de.xlf")

Write-Host "Handback report generated."
